$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the coordinate values in row 13 (point index 8) - B/C/D columns.
$ws.Range("B13").Value = 10561
$ws.Range("C13").Value = 20561
$ws.Range("D13").Value = 30561

# Update the active selection to match the saved cursor position (C17).
$ws.Range("C17").Select()
